# Sciwheel API - About page: added to db and increased randomness
#
# This script reproduces, via the Excel COM object model, the edits made to
# "New Features implementation plan-2020.xlsx":
#   1. The "Time spent" column of Table1 was right-aligned.
#   2. A handful of task names were edited / clarified.
#   3. The "Time spent" entries for two tasks were converted from free text
#      ("1hr20mins") or refined (3 -> 3.25 hours) into numeric hour values.
#   4. A new task row ("Write copy for achievements") was highlighted the
#      same way as the existing "Decide on content and website placement"
#      row, which also made that row wrap onto two lines.
#   5. The selection was left on B17.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Right-align the whole "Time spent" table column (data cells only,
#    header keeps its default formatting).
# ---------------------------------------------------------------------
$table = $ws.ListObjects.Item(1)
$timeSpentColumn = $table.ListColumns.Item(4)
$timeSpentColumn.DataBodyRange.HorizontalAlignment = -4152   # xlRight

# ---------------------------------------------------------------------
# 2. Rename a few tasks. New shared-string entries get appended in the
#    order the cells are written, so keep this order: B18, B19, B22.
# ---------------------------------------------------------------------
$ws.Range("B18").Value = "Research display format"

# "Write copy" -> "Write copy for achievements": also pick up the same
# highlight formatting used on the "Decide on content and website
# placement" task (row 17), and grow the row since the text now wraps.
$ws.Range("B17").Copy() | Out-Null
$ws.Range("B19").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Range("B19").Value = "Write copy for achievements"
$ws.Rows.Item(19).RowHeight = 29

$ws.Range("B22").Value = "Design list front-end"

# ---------------------------------------------------------------------
# 3. Update the "Time spent" numbers.
# ---------------------------------------------------------------------
$ws.Range("D17").Value = 2.5      # was free text "1hr20mins"
$ws.Range("D21").Value = 3.25     # was 3

# ---------------------------------------------------------------------
# 4. Leave the selection on B17, matching the saved view.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("B17").Select() | Out-Null

$wb.Save()
